# Products sample workbook edit:
#   - Drop the "CurrencyId" column entirely (it shifts NetPrice / NetShippingPrice /
#     Mrp / Comments one column to the left).
#   - Drop the trailing "ManufacturerPartNumber", "Gtin" and "Status" columns
#     (their header text + shared strings disappear; the now-unused column to
#     their right collapses out of the sheet's used range).
#   - Reset the view back to the top-left of the sheet (the old view had
#     scrolled all the way out to the now-deleted columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "CurrencyId" column (column Q) -- shifts everything after it left.
$ws.Range("Q1").EntireColumn.Delete()

# Remove the trailing "ManufacturerPartNumber", "Gtin" and "Status" columns
# (now columns U, V, W after the shift above).
$ws.Range("U1:W2").ClearContents()

# Bring the view/selection back on screen instead of the old far-right cell.
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
